$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2249443616207714
$ws.Range("C2").Value = 0.6609492819494724
$ws.Range("D2").Value = 0.739933022762375
$ws.Range("E2").Value = 0.8601935960947251
$ws.Range("F2").Value = 0.8616023889833488

$ws.Range("B3").Value = 0.0646333293875067
$ws.Range("C3").Value = 0.6340343087035507
$ws.Range("D3").Value = 0.9392448356363884
$ws.Range("E3").Value = 0.9691464469502988
$ws.Range("F3").Value = 1.019295709120695

$ws.Range("B4").Value = 0.47451964639446
$ws.Range("C4").Value = 0.47451964639446
$ws.Range("D4").Value = 0.4501379313770976
$ws.Range("E4").Value = 0.6709231933515919
$ws.Range("F4").Value = 0.5195794875428871

$ws.Range("B5").Value = 0.2995951651330526
$ws.Range("C5").Value = 0.2995951651330526
$ws.Range("D5").Value = 0.08977393156539534
$ws.Range("E5").Value = 0.299622982371839
$ws.Range("F5").Value = 0.005773836557138828
